# edit in students import and second edit in it && edit in teacher profile
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Students import sheet (Sheet1) edits ---

# Insert a new column before S ("current_year") to make room for the new
# "nationality" column; this shifts the old S:V (current_year .. parent_name)
# one column to the right, to T:W, carrying their widths/formatting along.
$ws.Columns("S:S").Insert()

# New header cells
$ws.Range("S1").Value = "nationality"
$ws.Range("X1").Value = "secondary_parent_id"
$ws.Range("Y1").Value = "secondary_parent_name"

# Widen the two new trailing columns to match the template layout
$ws.Columns("X:X").ColumnWidth = 19.166666666666668
$ws.Columns("Y:Y").ColumnWidth = 21.5

# Row 2 (sample record) value updates
$ws.Range("A2").Value = "شهاب"
$ws.Range("B2").Value = "new_1"
$ws.Range("C2").Value = 5555512345
$ws.Range("D2").Value = "ice_magic@gogo.com"
$ws.Range("M2").Value = "shehab"
$ws.Range("N2").Value = "ahmed"
$ws.Range("O2").Value = "awd"
$ws.Range("S2").Value = "SA"
$ws.Range("V2").Value = 6644885522
$ws.Range("W2").Value = "new one"
$ws.Range("X2").Value = 5356644425
$ws.Range("Y2").Value = "khaled hussain"

# Style the email cell as a hyperlink-styled cell (adds the Hyperlink
# font/cellStyle, matching the template's "teacher profile" email styling)
$ws.Range("D2").Style = "Hyperlink"

# Restore selection to the new last data cell
$ws.Range("AA2").Select()
